$d = $word.ActiveDocument

# The document currently has one paragraph (Iván's note). Add a new
# paragraph after it for Ismael's note, matching the "registro de
# modificaciones" convention already used by the first entry.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$endRange = $lastPara.Range
$endRange.Collapse(0)  # wdCollapseEnd
$endRange.InsertParagraphAfter()

# Move into the newly created (still empty) paragraph and type its text.
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newRange = $newPara.Range
$newRange.Collapse(1)  # wdCollapseStart
$newRange.InsertAfter("Ismael: Agregado una función que no hace nada")

# Word stamps the site of the last edit with a hidden "_GoBack" bookmark,
# collapsed right at the end of the text we just typed (before the
# paragraph mark). Adding a bookmark exactly at a paragraph's
# end-of-text position can misplace it, so nudge past the boundary with
# a throwaway character, anchor the bookmark there, then remove it.
$textEndPos = $d.Content.End - 1
$guardRange = $d.Range($textEndPos, $textEndPos)
$guardRange.InsertAfter("X")

$bookmarkRange = $d.Range($textEndPos, $textEndPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$guardCharRange = $d.Range($textEndPos, $textEndPos + 1)
$guardCharRange.Delete()
